# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock) sheet gets three new trailing columns - date,
# legislator_name and legislator_id - populated on every existing data
# row with the filing date (2013-12-24), the legislator's name (蔡其昌)
# and numeric id (1377).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Match the existing header formatting (bold font + border) used by the
# other columns before writing the new header labels.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

for ($r = 2; $r -le $lastRow; $r++) {
    # Force the date column to plain text so "2013-12-24" is kept as a
    # literal string instead of being reinterpreted as a date serial
    # number, then drop back to the sheet's normal (default) formatting.
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = "2013-12-24"
    $ws.Cells.Item($r, 8).ClearFormats()

    $ws.Cells.Item($r, 9).Value = "蔡其昌"
    $ws.Cells.Item($r, 10).Value = 1377
}
